# Insert two new weekly records (rows 170 and 171) into the "Zapallo"
# (Hortaliza, Macroferia Regional de Talca) price sheet. All existing
# rows from the old row 170 onward shift down by two rows (old row 170
# becomes row 172, old row 279 becomes row 281, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 170.. down by two, inserting two blank rows at 170-171.
$ws.Rows("170:171").Insert()

# --- New row 170: Camote, "1a (guarda)", origin Region del Maule ---
$ws.Cells.Item(170, 1).Value  = 5
$ws.Cells.Item(170, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(170, 3).Value  = "Maule"
$ws.Cells.Item(170, 4).Value  = 44777
$ws.Cells.Item(170, 5).Value  = 7
$ws.Cells.Item(170, 6).Value  = 100112045
$ws.Cells.Item(170, 7).Value  = "Zapallo"
$ws.Cells.Item(170, 8).Value  = "Camote"
$ws.Cells.Item(170, 9).Value  = "1a (guarda)"
$ws.Cells.Item(170, 10).Value = 900
$ws.Cells.Item(170, 11).Value = 600
$ws.Cells.Item(170, 12).Value = 600
$ws.Cells.Item(170, 13).Value = 600
$ws.Cells.Item(170, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(170, 15).Value = "Región del Maule"
$ws.Cells.Item(170, 16).Value = 600
$ws.Cells.Item(170, 17).Value = 1
$ws.Cells.Item(170, 18).Value = "Hortaliza"

# --- New row 171: Paine, "1a (guarda)", origin Region del Maule ---
$ws.Cells.Item(171, 1).Value  = 5
$ws.Cells.Item(171, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(171, 3).Value  = "Maule"
$ws.Cells.Item(171, 4).Value  = 44777
$ws.Cells.Item(171, 5).Value  = 7
$ws.Cells.Item(171, 6).Value  = 100112045
$ws.Cells.Item(171, 7).Value  = "Zapallo"
$ws.Cells.Item(171, 8).Value  = "Paine"
$ws.Cells.Item(171, 9).Value  = "1a (guarda)"
$ws.Cells.Item(171, 10).Value = 1500
$ws.Cells.Item(171, 11).Value = 200
$ws.Cells.Item(171, 12).Value = 200
$ws.Cells.Item(171, 13).Value = 200
$ws.Cells.Item(171, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(171, 15).Value = "Región del Maule"
$ws.Cells.Item(171, 16).Value = 200
$ws.Cells.Item(171, 17).Value = 1
$ws.Cells.Item(171, 18).Value = "Hortaliza"
